# Generate Report for Handback
# The e2e/ebb52c47-83da-4006-ae86-ec3795f654d8.md file has completed its
# handback cycle: its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the stale "version mismatch" error
# clears, and each locale's "Latest Handback DateTime" is refreshed.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the ebb52c47... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet: row 3 is the ebb52c47... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("K3").Value = "2016-08-22 20:47:58"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the ebb52c47... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("K3").Value = "2016-08-22 20:48:15"
$wsDeDe.Range("P3").Value = ""

# Error Detail column no longer holds long text; let it shrink back down.
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
